$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.491.01"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.109.90"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5234"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4493"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").Value = "2.107.64"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.789"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001126"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06616"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "30.528.66"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("E25").Value = "  +3.31%  "

$ws.Range("D26").Value = "2.346.77"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.76%  "

$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.167"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E35").Value = "  -2.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02575"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06806"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.500"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2284"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6941"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.338"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.81%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6389"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.645"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("E50").Value = "  +5.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
